{"js": "// Replace each old equation string with its new value, cell by cell,\n// using Body.search() (exact text, case-sensitive) + Range.insertText(..., \"Replace\").\nconst replacements = [\n  [\"90-8=82\", \"97-89=8\"],\n  [\"94-50=44\", \"38+60=98\"],\n  [\"76-40=36\", \"84-10=74\"],\n  [\"4+54=58\", \"22-1=21\"],\n  [\"88-53=35\", \"50-36=14\"],\n  [\"8+37=45\", \"67+32=99\"],\n  [\"51+7=58\", \"20-1=19\"],\n  [\"7+17=24\", \"31+5=36\"],\n  [\"35-9=26\", \"78+0=78\"],\n  [\"46+53=99\", \"71-28=43\"],\n  [\"50+42=92\", \"11+10=21\"],\n  [\"68-62=6\", \"47-40=7\"],\n  [\"2+44=46\", \"47+15=62\"],\n  [\"20+9=29\", \"34+63=97\"],\n  [\"53+13=66\", \"59-19=40\"],\n  [\"52+12=64\", \"15+79=94\"],\n  [\"72-51=21\", \"30+38=68\"],\n  [\"70-45=25\", \"3+95=98\"],\n  [\"29-20=9\", \"26-2=24\"],\n  [\"19+72=91\", \"42-14=28\"],\n  [\"30-23=7\", \"81-46=35\"],\n  [\"22-19=3\", \"82-39=43\"],\n  [\"28+59=87\", \"10-3=7\"],\n  [\"55+36=91\", \"91-5=86\"],\n  [\"10+39=49\", \"57+17=74\"],\n  [\"61-46=15\", \"87-63=24\"],\n  [\"24+42=66\", \"26+24=50\"],\n  [\"44+12=56\", \"97-2=95\"],\n  [\"70-63=7\", \"25+60=85\"],\n  [\"46+4=50\", \"86-16=70\"],\n  [\"13+9=22\", \"77-43=34\"],\n  [\"31+62=93\", \"87-21=66\"],\n  [\"91-41=50\", \"88-28=60\"],\n  [\"94-68=26\", \"32-28=4\"],\n  [\"35+53=88\", \"13-1=12\"],\n  [\"52+18=70\", \"49+21=70\"],\n  [\"98-43=55\", \"55-43=12\"],\n  [\"33+29=62\", \"34-0=34\"],\n  [\"61-31=30\", \"96-5=91\"],\n  [\"92-48=44\", \"72-49=23\"],\n  [\"71+16=87\", \"46+11=57\"],\n  [\"28+25=53\", \"20-0=20\"],\n  [\"27+9=36\", \"4+29=33\"],\n  [\"67-59=8\", \"20+78=98\"],\n  [\"0+42=42\", \"33-3=30\"],\n  [\"48+4=52\", \"97-10=87\"],\n  [\"39+47=86\", \"89-7=82\"],\n  [\"50-0=50\", \"59-35=24\"],\n  [\"63+33=96\", \"42+17=59\"],\n  [\"23+8=31\", \"51+35=86\"],\n  [\"15+29=44\", \"56-2=54\"],\n  [\"50-6=44\", \"26+18=44\"],\n  [\"48-35=13\", \"28+32=60\"],\n  [\"88+3=91\", \"94-69=25\"],\n  [\"43+53=96\", \"5+62=67\"],\n  [\"23+60=83\", \"38+11=49\"],\n  [\"29-7=22\", \"53+30=83\"],\n  [\"95-23=72\", \"41+39=80\"],\n  [\"41+9=50\", \"22+17=39\"],\n  [\"13+1=14\", \"74-63=11\"],\n  [\"24-21=3\", \"61+20=81\"],\n  [\"66-29=37\", \"86-25=61\"],\n  [\"11+46=57\", \"38+37=75\"],\n  [\"79-1=78\", \"22+34=56\"],\n  [\"68+27=95\", \"65+6=71\"],\n  [\"79+13=92\", \"42+9=51\"],\n  [\"3+26=29\", \"39-34=5\"],\n  [\"45-14=31\", \"19-14=5\"],\n  [\"77-17=60\", \"10+30=40\"],\n  [\"72+15=87\", \"53-17=36\"],\n  [\"98-0=98\", \"92-91=1\"],\n  [\"38+53=91\", \"30+16=46\"],\n  [\"97-57=40\", \"4+80=84\"],\n  [\"77-59=18\", \"55-36=19\"],\n  [\"85-70=15\", \"10+57=67\"],\n  [\"11+9=20\", \"90-67=23\"],\n  [\"10+58=68\", \"32+28=60\"],\n  [\"12+45=57\", \"31+26=57\"],\n  [\"87-24=63\", \"82-58=24\"],\n  [\"64-22=42\", \"69-55=14\"],\n  [\"95-89=6\", \"38-13=25\"],\n  [\"38-8=30\", \"37+10=47\"],\n  [\"62+16=78\", \"35+27=62\"],\n  [\"86-41=45\", \"32-30=2\"],\n  [\"78-71=7\", \"16+60=76\"],\n  [\"23+62=85\", \"15+52=67\"],\n  [\"62-2=60\", \"61-27=34\"],\n  [\"4+15=19\", \"76-38=38\"],\n  [\"73-24=49\", \"82-54=28\"],\n  [\"32+33=65\", \"26+37=63\"],\n  [\"92-5=87\", \"8+61=69\"],\n  [\"83-10=73\", \"92-84=8\"],\n  [\"14+29=43\", \"45-22=23\"],\n  [\"80-79=1\", \"47+34=81\"],\n  [\"43+17=60\", \"37-7=30\"],\n  [\"3+15=18\", \"41+15=56\"],\n  [\"76-44=32\", \"33+34=67\"],\n  [\"12+86=98\", \"85-47=38\"],\n  [\"46+24=70\", \"20+70=90\"],\n  [\"33+28=61\", \"21+53=74\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each old equation string with its new value, cell by cell, using\n# Range.Find.Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards,\n# MatchSoundsLike, MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace)\n# with Replace = wdReplaceAll (2), scoped to the whole document body each time.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    ,@('90-8=82', '97-89=8')\n    ,@('94-50=44', '38+60=98')\n    ,@('76-40=36', '84-10=74')\n    ,@('4+54=58', '22-1=21')\n    ,@('88-53=35', '50-36=14')\n    ,@('8+37=45', '67+32=99')\n    ,@('51+7=58', '20-1=19')\n    ,@('7+17=24', '31+5=36')\n    ,@('35-9=26', '78+0=78')\n    ,@('46+53=99', '71-28=43')\n    ,@('50+42=92', '11+10=21')\n    ,@('68-62=6', '47-40=7')\n    ,@('2+44=46', '47+15=62')\n    ,@('20+9=29', '34+63=97')\n    ,@('53+13=66', '59-19=40')\n    ,@('52+12=64', '15+79=94')\n    ,@('72-51=21', '30+38=68')\n    ,@('70-45=25', '3+95=98')\n    ,@('29-20=9', '26-2=24')\n    ,@('19+72=91', '42-14=28')\n    ,@('30-23=7', '81-46=35')\n    ,@('22-19=3', '82-39=43')\n    ,@('28+59=87', '10-3=7')\n    ,@('55+36=91', '91-5=86')\n    ,@('10+39=49', '57+17=74')\n    ,@('61-46=15', '87-63=24')\n    ,@('24+42=66', '26+24=50')\n    ,@('44+12=56', '97-2=95')\n    ,@('70-63=7', '25+60=85')\n    ,@('46+4=50', '86-16=70')\n    ,@('13+9=22', '77-43=34')\n    ,@('31+62=93', '87-21=66')\n    ,@('91-41=50', '88-28=60')\n    ,@('94-68=26', '32-28=4')\n    ,@('35+53=88', '13-1=12')\n    ,@('52+18=70', '49+21=70')\n    ,@('98-43=55', '55-43=12')\n    ,@('33+29=62', '34-0=34')\n    ,@('61-31=30', '96-5=91')\n    ,@('92-48=44', '72-49=23')\n    ,@('71+16=87', '46+11=57')\n    ,@('28+25=53', '20-0=20')\n    ,@('27+9=36', '4+29=33')\n    ,@('67-59=8', '20+78=98')\n    ,@('0+42=42', '33-3=30')\n    ,@('48+4=52', '97-10=87')\n    ,@('39+47=86', '89-7=82')\n    ,@('50-0=50', '59-35=24')\n    ,@('63+33=96', '42+17=59')\n    ,@('23+8=31', '51+35=86')\n    ,@('15+29=44', '56-2=54')\n    ,@('50-6=44', '26+18=44')\n    ,@('48-35=13', '28+32=60')\n    ,@('88+3=91', '94-69=25')\n    ,@('43+53=96', '5+62=67')\n    ,@('23+60=83', '38+11=49')\n    ,@('29-7=22', '53+30=83')\n    ,@('95-23=72', '41+39=80')\n    ,@('41+9=50', '22+17=39')\n    ,@('13+1=14', '74-63=11')\n    ,@('24-21=3', '61+20=81')\n    ,@('66-29=37', '86-25=61')\n    ,@('11+46=57', '38+37=75')\n    ,@('79-1=78', '22+34=56')\n    ,@('68+27=95', '65+6=71')\n    ,@('79+13=92', '42+9=51')\n    ,@('3+26=29', '39-34=5')\n    ,@('45-14=31', '19-14=5')\n    ,@('77-17=60', '10+30=40')\n    ,@('72+15=87', '53-17=36')\n    ,@('98-0=98', '92-91=1')\n    ,@('38+53=91', '30+16=46')\n    ,@('97-57=40', '4+80=84')\n    ,@('77-59=18', '55-36=19')\n    ,@('85-70=15', '10+57=67')\n    ,@('11+9=20', '90-67=23')\n    ,@('10+58=68', '32+28=60')\n    ,@('12+45=57', '31+26=57')\n    ,@('87-24=63', '82-58=24')\n    ,@('64-22=42', '69-55=14')\n    ,@('95-89=6', '38-13=25')\n    ,@('38-8=30', '37+10=47')\n    ,@('62+16=78', '35+27=62')\n    ,@('86-41=45', '32-30=2')\n    ,@('78-71=7', '16+60=76')\n    ,@('23+62=85', '15+52=67')\n    ,@('62-2=60', '61-27=34')\n    ,@('4+15=19', '76-38=38')\n    ,@('73-24=49', '82-54=28')\n    ,@('32+33=65', '26+37=63')\n    ,@('92-5=87', '8+61=69')\n    ,@('83-10=73', '92-84=8')\n    ,@('14+29=43', '45-22=23')\n    ,@('80-79=1', '47+34=81')\n    ,@('43+17=60', '37-7=30')\n    ,@('3+15=18', '41+15=56')\n    ,@('76-44=32', '33+34=67')\n    ,@('12+86=98', '85-47=38')\n    ,@('46+24=70', '20+70=90')\n    ,@('33+28=61', '21+53=74')\n)\n\n$missing = @()\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $found = $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        $missing += $oldText\n    }\n}\n\nif ($missing.Count -gt 0) {\n    throw \"No match found for: \" + ($missing -join \", \")\n}\n"}
